# Generate Report for Handoff
# Update the "7bc79fcf-000e-47b3-a00c-ff4582d3354f.md" row's status to
# "Ready for handoff" across the Overview/zh-cn/de-de sheets, and record
# the new handoff datetimes for zh-cn and de-de.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is the 7bc79fcf-... entry (B = zh-cn status, C = de-de status)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the 7bc79fcf-... entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-02-17 06:12:38"

# de-de sheet: row 3 is the 7bc79fcf-... entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-02-17 06:12:48"
